# Updates COVID-19 BD workbook with data for 2020-04-13 (serial 43934)
# and recomputed mortality-rate (Proporción) figures.

$wb = $excel.ActiveWorkbook

$wsProp  = $wb.Worksheets.Item("Proporción")
$wsMun   = $wb.Worksheets.Item("Municipios")
$wsConf  = $wb.Worksheets.Item("Confirmados")

# ---------------------------------------------------------------------------
# Sheet "Proporción": add row 11
# ---------------------------------------------------------------------------

# Carry the date / percentage number formats down from row 10 into row 11
$wsProp.Range("A10").Copy() | Out-Null
$wsProp.Range("A11").PasteSpecial(-4122) | Out-Null
$wsProp.Range("I10").Copy() | Out-Null
$wsProp.Range("I11").PasteSpecial(-4122) | Out-Null
$wsProp.Range("K10:P10").Copy() | Out-Null
$wsProp.Range("K11:P11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsProp.Range("A11").Value2 = 43934
$wsProp.Range("B11").Formula = "=C11+D11+E11+F11+G11"
$wsProp.Range("C11").Value2 = 6
$wsProp.Range("D11").Value2 = 59
$wsProp.Range("E11").Value2 = 5
$wsProp.Range("F11").Value2 = 6
$wsProp.Range("G11").Value2 = 12

$wsProp.Range("I11").Value2 = 43934
$wsProp.Range("J11").Formula = "=B11"
$wsProp.Range("K11").Formula = "=C11/`$B11"
$wsProp.Range("L11").Formula = "=D11/`$B11"
$wsProp.Range("M11").Formula = "=E11/`$B11"
$wsProp.Range("N11").Formula = "=F11/`$B11"
$wsProp.Range("O11").Formula = "=G11/`$B11"
$wsProp.Range("P11").Formula = "=SUM(K11:O11)"

# ---------------------------------------------------------------------------
# Sheet "Municipios": add row 11
# ---------------------------------------------------------------------------

$wsMun.Range("A10").Copy() | Out-Null
$wsMun.Range("A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsMun.Range("A11").Value2 = 43934
$wsMun.Range("B11").Formula = "=SUM(C11:M11)"
$wsMun.Range("C11").Value2 = 24
$wsMun.Range("D11").Value2 = 7
$wsMun.Range("E11").Value2 = 31
$wsMun.Range("F11").Value2 = 4
$wsMun.Range("G11").Value2 = 2
$wsMun.Range("H11").Value2 = 4
$wsMun.Range("I11").Value2 = 1
$wsMun.Range("J11").Value2 = 8
$wsMun.Range("K11").Value2 = 3
$wsMun.Range("L11").Value2 = 2
$wsMun.Range("M11").Value2 = 2

# ---------------------------------------------------------------------------
# Sheet "Confirmados": add row 31
# ---------------------------------------------------------------------------

$wsConf.Range("A30").Copy() | Out-Null
$wsConf.Range("A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsConf.Range("A31").Value2 = 43934
$wsConf.Range("B31").Value2 = 2
$wsConf.Range("C31").Formula = "=C30+B31"

# ---------------------------------------------------------------------------
# View state: update per-sheet selections, zoom on "Proporción", and make
# "Confirmados" the active/selected tab (matches the saved workbook state).
# ---------------------------------------------------------------------------

$wsProp.Activate()
$wsProp.Range("F11").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110

$wsMun.Activate()
$wsMun.Range("C11").Select() | Out-Null

$wsConf.Activate()
$wsConf.Range("D34").Select() | Out-Null

Write-Host "Workbook updated with 2020-04-13 data."
